# Update Fecha (D), Calidad (L), Volumen (M), Precio minimo (N),
# Precio maximo (O), Precio promedio ponderado (P) and Precio $/Kg (S)
# for rows 2-12 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param($Row, $Fecha, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $PrecioKg)

    $ws.Cells.Item($Row, 4).Value  = $Fecha       # D - Fecha
    $ws.Cells.Item($Row, 12).Value = $Calidad     # L - Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen     # M - Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin   # N - Precio minimo
    $ws.Cells.Item($Row, 15).Value = $PrecioMax   # O - Precio maximo
    $ws.Cells.Item($Row, 16).Value = $PrecioProm  # P - Precio promedio ponderado
    $ws.Cells.Item($Row, 19).Value = $PrecioKg    # S - Precio $/Kg
}

Set-Row 2  44435 "Primera"  40  20000 20000 20000 2000
Set-Row 3  44432 "Primera"  20  20000 20000 20000 2000
Set-Row 4  44466 "Primera"  60  20000 20000 20000 2000
Set-Row 5  44511 "Primera"  120 28000 28000 28000 2800
Set-Row 6  44473 "Primera"  180 20000 20000 20000 2000
Set-Row 7  44503 "Primera"  60  30000 30000 30000 3000
Set-Row 8  44503 "Segunda"  50  25000 25000 25000 2500
Set-Row 9  44476 "Primera"  120 20000 20000 20000 2000
Set-Row 10 44434 "Primera"  20  20000 20000 20000 2000
Set-Row 11 44517 "Especial" 100 27000 27000 27000 2700
Set-Row 12 44517 "Primera"  30  25000 25000 25000 2500
